# Apply the "text" function-group addition to the hidden '#system' sheet,
# mirroring the commit:
# [base] - [`outputToCloud(resource)`]: support the transferring of output
# artifact to the cloud.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Column A ("target" list of group names): insert a new row for the
#    "text" group right before "web" (row 25), shifting web..xml down by
#    one row (25-30 -> 26-31).
# ---------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $v = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 1).Value = $v
}
$ws.Cells.Item(25, 1).Value = "text"

# ---------------------------------------------------------------------
# 2) Column E ("base" functions): insert "outputToCloud(resource)" right
#    before "prependText(var,prependWith)" (row 22), shifting the rest
#    of the alphabetical list down by one row (22-38 -> 23-39).
# ---------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $v = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r + 1, 5).Value = $v
}
$ws.Cells.Item(22, 5).Value = "outputToCloud(resource)"

# ---------------------------------------------------------------------
# 3) Insert a whole new column at Y so that the existing web / webalert /
#    webcookie / ws / ws.async / xml columns (Y..AD) shift right to
#    Z..AE, then populate the freed column Y with the new "text" group
#    (header + single function entry).
# ---------------------------------------------------------------------
$ws.Columns("Y:Y").Insert()
$ws.Cells.Item(1, 25).Value = "text"
$ws.Cells.Item(2, 25).Value = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------
# 4) Fix up the defined names so they point at the correct (shifted)
#    ranges, and register the brand-new "text" name.
# ---------------------------------------------------------------------
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"

$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")

Write-Host "Applied 'text' group / outputToCloud(resource) edits."
